$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.984.54"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").Value = "1.560.33"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.490"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.91%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("E8").Value = "  +2.03%  "

$ws.Range("E9").Value = "  +0.16%  "

$ws.Range("E10").Value = "  +1.85%  "

$ws.Range("E11").Value = "  +0.28%  "

$ws.Range("D12").Value = "1.782.03"

$ws.Range("D13").Value = "1.545.26"
$ws.Range("E13").Value = "  -0.49%  "

$ws.Range("E15").Value = "  +1.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.69%  "

$ws.Range("D17").Value = "26.977.96"
$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").Value = "0.0₃0705"
$ws.Range("E18").Value = "  +2.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.42%  "

$ws.Range("E20").Value = "  +2.34%  "

$ws.Range("E21").Value = "  -0.17%  "

$ws.Range("E22").Value = "  +1.37%  "

$ws.Range("E23").Value = "  +0.50%  "

$ws.Range("E24").Value = "  -1.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.05"
$ws.Range("D27").Style = "Normal"

$ws.Range("E28").Value = "  +1.42%  "

$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("E30").Value = "  +0.74%  "

$ws.Range("E31").Value = "  +1.99%  "

$ws.Range("E32").Value = "  +0.70%  "

$ws.Range("D33").Value = "1.424.93"
$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("E34").Value = "  +3.77%  "

$ws.Range("E35").Value = "  +3.05%  "

$ws.Range("E36").Value = "  +9.68%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0165"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.62%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.530"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.808"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.16%  "

$ws.Range("E41").Value = "  -0.14%  "

$ws.Range("E42").Value = "  +0.63%  "

$ws.Range("E43").Value = "  +2.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.89%  "

$ws.Range("E46").Value = "  -0.20%  "

$ws.Range("D47").Value = "1.695.84"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.32%  "

$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("D50").Value = "0.0₆01000"
$ws.Range("E50").Value = "  -1.10%  "

$ws.Range("E51").Value = "  -0.11%  "

